$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Fix the typo in "Make a new bgitfor the release: " so it reads
# "Make a new branch for the release: ", and leave Word's "_GoBack"
# last-edit bookmark positioned right after the newly typed "ranch ".
# ------------------------------------------------------------------

# Locate the "b" / "git" / "for" runs that spell out the typo "bgitfor".
$rFind = $d.Content
$found = $rFind.Find.Execute("bgitfor", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$bStart   = $rFind.Start   # start of "b"
$gitStart = $bStart + 1    # start of "git"   ("b" is always 1 char)
$forStart = $gitStart + 3  # start of "for"   ("git" is always 3 chars)

# Temporary marker between "b" and "git" - keeps the upcoming "git"->"ranch "
# replacement from being absorbed into the "b" run while we still need them
# to stay apart.
$d.Bookmarks.Add("zzTempSplit", $d.Range($gitStart, $gitStart))

# Word's "_GoBack" bookmark is a singleton: adding it here both creates it
# at the new location and removes it from wherever it used to be.
$d.Bookmarks.Add("_GoBack", $d.Range($forStart, $forStart))

# Merge "Make a new " and "b" into a single run ("Make a new b"). Setting a
# range's .Text to a value that actually differs, then immediately back to
# the original text, forces the two adjacent runs to coalesce.
$rMergeB = $d.Range($bStart - 1, $bStart + 1)
$originalMergeB = $rMergeB.Text
$rMergeB.Text = "#"
$d.Range($bStart - 1, $bStart).Text = $originalMergeB

# Replace "git" with "ranch " - becomes its own run because it is fenced in
# by the temp marker on the left and the _GoBack bookmark on the right.
$d.Range($gitStart, $forStart).Text = "ranch "

# The temp marker has done its job; removing it does not re-merge runs.
$d.Bookmarks("zzTempSplit").Delete()

# Merge "for" and " the release: " into a single run ("for the release: ").
# "ranch " now ends right where the _GoBack bookmark sits, so this edit
# (just past the bookmark) cannot reach back into it.
$full = $d.Content.Text
$idx = $full.IndexOf("for the release")
$rMergeFor = $d.Range($idx + 2, $idx + 4)
$originalMergeFor = $rMergeFor.Text
$rMergeFor.Text = "#X"
$d.Range($idx + 2, $idx + 4).Text = $originalMergeFor
